$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, which shifts all existing rows (4..111) down
# by one (to 5..112), matching the rest of the diff automatically.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record (week's new entry).
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13425
$ws.Range("N4").Value = "$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1033
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
